# Update of league bases, 17-06-2024 21:10
# The source feed reshuffled which match record lands on which row for a
# few fixtures (ids moved between rows while the row's "id" index column A
# stayed put). Reproduce that by moving the B:AD payload between rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Swap-Rows($row1, $row2) {
    $r1 = $ws.Range("B$row1`:AD$row1")
    $r2 = $ws.Range("B$row2`:AD$row2")
    $v1 = $r1.Value()
    $v2 = $r2.Value()
    $r1.Value = $v2
    $r2.Value = $v1
}

# Rows 130 and 131 (match ids 7483189 / 7483081) swap places.
Swap-Rows 130 131

# Rows 139 and 140 (match ids 7528859 / 7528849) swap places.
Swap-Rows 139 140

# Rows 254, 255, 256 rotate: new254 = old255, new255 = old256, new256 = old254.
$v254 = $ws.Range("B254:AD254").Value()
$v255 = $ws.Range("B255:AD255").Value()
$v256 = $ws.Range("B256:AD256").Value()

$ws.Range("B254:AD254").Value = $v255
$ws.Range("B255:AD255").Value = $v256
$ws.Range("B256:AD256").Value = $v254
